{"js": "// Expand \"VSpy\" -> \"Vehicle Spy\" and append \", Jira, HIL\" after \"8D\" in the\n// SKILLS / Tools: line, as described by the commit\n// \"added some other skills for amphenol etc.\"\n\n// --- Edit 1: \"Intrepid VSpy, Vector CANoe,\" -> \"Intrepid Vehicle Spy, Vector CANoe,\" ---\n// Search for the unique tail \"Spy, Vector CANoe,\" and insert \"ehicle \" right\n// before it -- this lands the new text between the \"V\" and \"Spy\" of \"VSpy\",\n// turning it into \"Vehicle Spy\" while leaving all surrounding text untouched.\nconst vspyResults = context.document.body.search(\"Spy, Vector CANoe,\", { matchCase: true });\nvspyResults.load(\"items\");\nawait context.sync();\n\nif (vspyResults.items.length === 0) {\n  throw new Error('Could not find \"Spy, Vector CANoe,\" in the document body.');\n}\nvspyResults.items[0].insertText(\"ehicle \", \"Before\");\nawait context.sync();\n\n// --- Edit 2: \" ... Five Why, 8D\" -> \" ... Five Why, 8D, Jira, HIL\" ---\n// Search for the unique trailing token \"Five Why, 8D\" and insert the new\n// skills right after it.\nconst toolsResults = context.document.body.search(\"Five Why, 8D\", { matchCase: true });\ntoolsResults.load(\"items\");\nawait context.sync();\n\nif (toolsResults.items.length === 0) {\n  throw new Error('Could not find \"Five Why, 8D\" in the document body.');\n}\ntoolsResults.items[0].insertText(\", Jira, HIL\", \"After\");\nawait context.sync();\n", "ps1": "# Expand \"VSpy\" -> \"Vehicle Spy\" and append \", Jira, HIL\" after \"8D\" in the\n# SKILLS / Tools: line, as described by the commit\n# \"added some other skills for amphenol etc.\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"Intrepid VSpy, Vector CANoe,\" -> \"Intrepid Vehicle Spy, Vector CANoe,\" ---\n# Find the unique tail \"Spy, Vector CANoe,\", collapse the found range to its\n# start (a zero-width insertion point sitting right after the \"V\" of \"VSpy\"),\n# and insert \"ehicle \" there so \"VSpy\" becomes \"Vehicle Spy\".\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"Spy, Vector CANoe,\")\nif ($found1) {\n  $rng1.Collapse(1)\n  $rng1.Text = \"ehicle \"\n} else {\n  throw \"Could not find 'Spy, Vector CANoe,' in the document.\"\n}\n\n# --- Edit 2: \" ... Five Why, 8D\" -> \" ... Five Why, 8D, Jira, HIL\" ---\n# Find the unique trailing token \"Five Why, 8D\", collapse to its end, and\n# insert the new skills there.\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"Five Why, 8D\")\nif ($found2) {\n  $rng2.Collapse(0)\n  $rng2.Text = \", Jira, HIL\"\n} else {\n  throw \"Could not find 'Five Why, 8D' in the document.\"\n}\n"}
